# "Generate Report for Handback" - populate the "Latest Target File" /
# "Latest Handback File" columns (F/G) for the zh-cn and de-de handback
# report sheets, update the handback status + timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status -> "Handed back: in sync with en-US" (shared by both data rows)
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime
$wsZh.Range("H2").Value = "2016-03-20 08:37:28"
$wsZh.Range("H3").Value = "2016-03-20 08:37:28"

# Row 2 (79d2a641-...): Latest Target File / Latest Handback File
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ddb8e3c7fbe6aaeb688f4325d578aedd1a51015/e2e/79d2a641-b11c-4d0e-b0b3-cee54dd48687.md", "", "", "79d2a641-b11c-4d0e-b0b3-cee54dd48687.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/26b2cc1434141a7f1db44cbda039184af6a4510c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/79d2a641-b11c-4d0e-b0b3-cee54dd48687.7c4116a52db88b6ac2402d952247cb9e1fa7bdf3.zh-cn.xlf", "", "", "79d2a641-b11c-4d0e-b0b3-cee54dd48687.7c4116a52db88b6ac2402d952247cb9e1fa7bdf3.zh-cn.xlf") | Out-Null

# Row 3 (8c88d43a-...): Latest Target File / Latest Handback File
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/6ddb8e3c7fbe6aaeb688f4325d578aedd1a51015/e2e/8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md", "", "", "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/26b2cc1434141a7f1db44cbda039184af6a4510c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8c88d43a-4afe-4bc6-8302-f7f8b4106faf.50c414ecd8e910df6e3226df088baf65c337638d.zh-cn.xlf", "", "", "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.50c414ecd8e910df6e3226df088baf65c337638d.zh-cn.xlf") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status -> "Handed back: in sync with en-US" (shared by both data rows)
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# Latest Handback DateTime (new timestamp for de-de handback)
$wsDe.Range("H2").Value = "2016-03-20 08:37:34"
$wsDe.Range("H3").Value = "2016-03-20 08:37:34"

# Row 2 (79d2a641-...): Latest Target File / Latest Handback File
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/6ddb8e3c7fbe6aaeb688f4325d578aedd1a51015/e2e/79d2a641-b11c-4d0e-b0b3-cee54dd48687.md", "", "", "79d2a641-b11c-4d0e-b0b3-cee54dd48687.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4ed9f3ed49f967306434203f9472a1320524a7ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/79d2a641-b11c-4d0e-b0b3-cee54dd48687.7c4116a52db88b6ac2402d952247cb9e1fa7bdf3.de-de.xlf", "", "", "79d2a641-b11c-4d0e-b0b3-cee54dd48687.7c4116a52db88b6ac2402d952247cb9e1fa7bdf3.de-de.xlf") | Out-Null

# Row 3 (8c88d43a-...): Latest Target File / Latest Handback File
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/6ddb8e3c7fbe6aaeb688f4325d578aedd1a51015/e2e/8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md", "", "", "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4ed9f3ed49f967306434203f9472a1320524a7ac/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8c88d43a-4afe-4bc6-8302-f7f8b4106faf.50c414ecd8e910df6e3226df088baf65c337638d.de-de.xlf", "", "", "8c88d43a-4afe-4bc6-8302-f7f8b4106faf.50c414ecd8e910df6e3226df088baf65c337638d.de-de.xlf") | Out-Null
